# "cleaned up some whitespace" — the only substantive content change in the
# commit is on the slide with SlideID 277 ("Making it tangible" / "See it in
# real code"): the second paragraph's single run is split into two runs so
# the text reads "See working code" instead of "See it in real code".
#
# (The underlying OOXML diff also touches collaboration/revision bookkeeping
# parts — revisionInfo.xml, changesInfos/changesInfo1.xml — and the cached
# datetimeFigureOut text in the handout/notes masters. Those are
# PowerPoint-internal, session-generated metadata that aren't reachable
# through the public PowerPoint object model, so this script focuses on the
# real, user-visible edit.)

$p = $ppt.ActivePresentation

# Locate the slide by its stable SlideID (277) rather than a raw index, so
# the edit keeps working even if slide order changes.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 277) {
        $targetSlide = $candidate
    }
}
if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.Item(6)
}

# Locate the "Text Placeholder 3" shape (Id=4) holding the two lines of text.
$targetShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $candidate = $targetSlide.Shapes.Item($i)
    if ($candidate.Id -eq 4) {
        $targetShape = $candidate
    }
}
if ($targetShape -eq $null) {
    $targetShape = $targetSlide.Shapes.Item(1)
}

$textRange = $targetShape.TextFrame.TextRange
$paragraph = $textRange.Paragraphs(2)

# "See it in real code" -> "See working code", keeping "code" as its own
# trailing run (matching the two <a:r> runs in the target markup).
$paragraph.Runs(1).Text = "See working code"
$paragraph = $textRange.Paragraphs(2)
$len = $paragraph.Text.Length
$codeRun = $paragraph.Characters($len - 3, 4)
$codeRun.Text = "code"
